$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 3 values
$ws.Range("B3").Value = 1000
$ws.Range("D3").Value = 50

# Add new rows of data (12-15)
$data = @(
    @("bife do vazio", 1000, "g", 60),
    @("alho poro", 100, "Un", 5),
    @("cebola roxa", 500, "g", 7),
    @("arroz parboilizado", 1000, "Un", 30)
)

$row = 12
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row++
}
